$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) store numeric-looking text as literal
# strings (inlineStr) in the source workbook, e.g. "245.00" or "-0.35%".
# Force Text formatting while writing so Excel does not silently coerce
# these into numeric cells, then clear the format back off so no stray
# cell style is left behind on cells that originally had none.
$numericTextRange = $ws.Range("D2:E50")
$numericTextRange.NumberFormat = "@"

$ws.Range("D2").Value = '244.86'
$ws.Range("E2").Value = '-0.26%'
$ws.Range("E3").Value = '4.68%'
$ws.Range("D4").Value = '5.126'
$ws.Range("E4").Value = '0.39%'
$ws.Range("D5").Value = '0.05587'
$ws.Range("E5").Value = '0.19%'
$ws.Range("D6").Value = '6.467'
$ws.Range("E6").Value = '-0.49%'
$ws.Range("D7").Value = '0.8177'
$ws.Range("E7").Value = '0.04%'
$ws.Range("D8").Value = '0.8360'
$ws.Range("E8").Value = '-1.11%'
$ws.Range("B9").Value = 'One'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D9").Value = '0.0006002'
$ws.Range("E9").Value = '-93.88%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1333'
$ws.Range("E10").Value = '-0.63%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.06978'
$ws.Range("E11").Value = '0.37%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '0.02888'
$ws.Range("E12").Value = '-0.11%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '0.09390'
$ws.Range("E13").Value = '0.03%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '0.001519'
$ws.Range("E14").Value = '-0.39%'
$ws.Range("D15").Value = '0.006159'
$ws.Range("E15").Value = '-0.56%'
$ws.Range("D16").Value = '3.649'
$ws.Range("E16").Value = '4.23%'
$ws.Range("D17").Value = '3.036'
$ws.Range("E17").Value = '0.52%'
$ws.Range("D18").Value = '2.183'
$ws.Range("E18").Value = '5.81%'
$ws.Range("E19").Value = '-2.12%'
$ws.Range("D20").Value = '0.03112'
$ws.Range("E20").Value = '-1.70%'
$ws.Range("E21").Value = '-2.25%'
$ws.Range("D22").Value = '3.753'
$ws.Range("E22").Value = '0.15%'
$ws.Range("D23").Value = '0.04634'
$ws.Range("E23").Value = '-1.87%'
$ws.Range("E24").Value = '-0.10%'
$ws.Range("D25").Value = '0.001249'
$ws.Range("E25").Value = '-0.09%'
$ws.Range("D26").Value = '0.004502'
$ws.Range("E26").Value = '-2.93%'
$ws.Range("D27").Value = '0.00009604'
$ws.Range("E27").Value = '-1.04%'
$ws.Range("D28").Value = '0.0001394'
$ws.Range("E28").Value = '0.38%'
$ws.Range("D40").Value = '0.03643'
$ws.Range("E40").Value = '-0.73%'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = '0.006203'
$ws.Range("E41").Value = '83.58%'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '0.1051'
$ws.Range("E42").Value = '-22.97%'
$ws.Range("D43").Value = '0.002401'
$ws.Range("E43").Value = '-9.10%'
$ws.Range("D44").Value = '0.008855'
$ws.Range("E44").Value = '6.60%'
$ws.Range("D45").Value = '0.00005358'
$ws.Range("E45").Value = '1.06%'
$ws.Range("E46").Value = '0.00%'
$ws.Range("E47").Value = '-4.00%'
$ws.Range("D48").Value = '0.002339'
$ws.Range("E48").Value = '10.19%'
$ws.Range("D49").Value = '0.00002101'
$ws.Range("E49").Value = '0.00%'
$ws.Range("D50").Value = '0.0002001'
$ws.Range("E50").Value = '0.00%'

$numericTextRange.ClearFormats()
